# Fix formatting on fastq purpose column: "fullRNASEQ" -> "fullRNASeq"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E = purpose
    if ($cell.Value2 -eq "fullRNASEQ") {
        $cell.Value2 = "fullRNASeq"
    }
}
